$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1023.5
$ws.Range("I19").Value = 899.5
$ws.Range("J19").Value = 1147.5
$ws.Range("K19").Value = 899.5
$ws.Range("L19").Value = 1147.5
$ws.Range("M19").Value = -724.5
$ws.Range("N19").Value = -1497.5

$ws.Range("H33").Value = 222.86667
$ws.Range("I33").Value = 222.86667
$ws.Range("K33").Value = 222.86667
$ws.Range("M33").Value = 6.133330000000001

$ws.Range("H80").Value = 384.0625
$ws.Range("I80").Value = 190.125
$ws.Range("K80").Value = 570.375
$ws.Range("M80").Value = 427.625

$ws.Range("H83").Value = 384.0625
$ws.Range("I83").Value = 190.125
$ws.Range("K83").Value = 1711.125
$ws.Range("M83").Value = 3280.875

$ws.Range("H113").Value = 4119.3335
$ws.Range("I113").Value = 2776.4
$ws.Range("K113").Value = 2776.4
$ws.Range("M113").Value = 477.5999999999999

$ws.Range("H116").Value = 11075.6
$ws.Range("J116").Value = 14630
$ws.Range("L116").Value = 14630
$ws.Range("N116").Value = -21514

$ws.Range("H132").Value = 6996.9
$ws.Range("I132").Value = 7621.6113
$ws.Range("K132").Value = 22864.8339
$ws.Range("M132").Value = -20334.8339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 34147.69
$ws.Range("J45").Value = 1867.6666
$ws.Range("L45").Value = 1867.6666
$ws.Range("N45").Value = -2621.6666

$ws.Range("H97").Value = 577.8
$ws.Range("I97").Value = 494.13333
$ws.Range("K97").Value = 494.13333
$ws.Range("M97").Value = 1.866669999999999

$ws.Range("H132").Value = 2196.5833
$ws.Range("I132").Value = 1748.1177
$ws.Range("J132").Value = 3285.7144
$ws.Range("K132").Value = 5244.3531
$ws.Range("L132").Value = 9857.143199999999
$ws.Range("M132").Value = -2714.3531
$ws.Range("N132").Value = -14917.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3500.8928
$ws.Range("I86").Value = 3410.3044
$ws.Range("K86").Value = 3410.3044
$ws.Range("M86").Value = -2287.3044

$ws.Range("H89").Value = 3500.8928
$ws.Range("I89").Value = 3410.3044
$ws.Range("K89").Value = 17051.522
$ws.Range("M89").Value = -11435.522

$ws.Range("H106").Value = 29101.285
$ws.Range("J106").Value = 29101.285
$ws.Range("L106").Value = 29101.285
$ws.Range("N106").Value = -31625.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 780.6
$ws.Range("I5").Value = 780.6
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 780.6
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -668.6
$ws.Range("N5").ClearContents()

$ws.Range("H31").Value = 5684870.5
$ws.Range("I31").Value = 2196.5881
$ws.Range("K31").Value = 2196.5881
$ws.Range("M31").Value = -1901.5881

$ws.Range("H34").Value = 5684870.5
$ws.Range("I34").Value = 2196.5881
$ws.Range("K34").Value = 2196.5881
$ws.Range("M34").Value = -1994.5881

$ws.Range("H99").Value = 3979.625
$ws.Range("I99").Value = 3399.7
$ws.Range("K99").Value = 3399.7
$ws.Range("M99").Value = -1901.7

$ws.Range("H122").Value = 2087.4814
$ws.Range("I122").Value = 1794.7727
$ws.Range("J122").Value = 3375.4
$ws.Range("K122").Value = 5384.3181
$ws.Range("L122").Value = 10126.2
$ws.Range("M122").Value = -2934.3181
$ws.Range("N122").Value = -15026.2

$ws.Range("H126").Value = 3979.625
$ws.Range("I126").Value = 3399.7
$ws.Range("K126").Value = 10199.1
$ws.Range("M126").Value = -7729.099999999999

$ws.Range("H132").Value = 4113.1562
$ws.Range("J132").Value = 5272.636
$ws.Range("L132").Value = 15817.908
$ws.Range("N132").Value = -20877.908

$ws.Range("H134").Value = 4356.9033
$ws.Range("J134").Value = 4044.4285
$ws.Range("L134").Value = 12133.2855
$ws.Range("N134").Value = -17203.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 7697715.5
$ws.Range("I68").Value = 969.1667
$ws.Range("K68").Value = 2907.5001
$ws.Range("M68").Value = -2096.5001

$ws.Range("H71").Value = 7697715.5
$ws.Range("I71").Value = 969.1667
$ws.Range("K71").Value = 8722.5003
$ws.Range("M71").Value = -4666.5003

$ws.Range("H133").Value = 3000
$ws.Range("I133").Value = 1000
$ws.Range("K133").Value = 3000
$ws.Range("M133").Value = 2060

$ws.Range("H134").Value = 4999.7144
$ws.Range("I134").Value = 4998
$ws.Range("K134").Value = 14994
$ws.Range("M134").Value = -9924

$ws.Range("H138").Value = 20000
$ws.Range("J138").Value = 20000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 58825830
$ws.Range("J80").Value = 2386.5557
$ws.Range("L80").Value = 2386.5557
$ws.Range("N80").Value = -4382.5557

$ws.Range("H83").Value = 58825830
$ws.Range("J83").Value = 2386.5557
$ws.Range("L83").Value = 11932.7785
$ws.Range("N83").Value = -21916.7785

$ws.Range("H102").Value = 4367.1714
$ws.Range("I102").Value = 1760.125
$ws.Range("J102").Value = 5139.6294
$ws.Range("K102").Value = 1760.125
$ws.Range("L102").Value = 5139.6294
$ws.Range("M102").Value = -138.125
$ws.Range("N102").Value = -8383.6294

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 2392.6155
$ws.Range("I132").Value = 2321.4736
$ws.Range("J132").Value = 2585.7144
$ws.Range("K132").Value = 6964.4208
$ws.Range("L132").Value = 7757.1432
$ws.Range("M132").Value = -4434.4208
$ws.Range("N132").Value = -12817.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5325.2607
$ws.Range("I136").Value = 4865.3887
$ws.Range("J136").Value = 6980.8
$ws.Range("K136").Value = 14596.1661
$ws.Range("L136").Value = 20942.4
$ws.Range("M136").Value = -12046.1661
$ws.Range("N136").Value = -26042.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()

$ws.Range("H29").Value = 30010
$ws.Range("I29").Value = 30010
$ws.Range("K29").Value = 30010
$ws.Range("M29").Value = -29720

$ws.Range("H107").Value = 544.3158
$ws.Range("I107").Value = 517.86664
$ws.Range("K107").Value = 1553.59992
$ws.Range("M107").Value = 366.4000800000001

$ws.Range("H126").Value = 10923.538
$ws.Range("I126").Value = 13000.7
$ws.Range("K126").Value = 39002.10000000001
$ws.Range("M126").Value = -36532.10000000001

